# lhs_settings_input.xlsx — CST EB updated run files
# Applies the data edits described in the commit diff:
#  - parameters sheet: new "u" distribution labels for rows 2-5, new
#    value-1/value-2 numbers for rows 2-5 and row 13 (D13)
#  - initial_conditions sheet: value-1/value-2 bumped from 2 to 20
#  - active sheet switches from "parameters" to "initial_conditions"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("parameters")
$ws2 = $wb.Worksheets.Item("initial_conditions")

# --- parameters sheet -------------------------------------------------

# Rows 2-4 (k_grow parameters): give B2/B3 the same "text" formatting that
# B4 already carries, then fill in the distribution label + new values.
$ws1.Range("B4").Copy()
$ws1.Range("B2").PasteSpecial(-4122)
$ws1.Range("B3").PasteSpecial(-4122)

$ws1.Range("B2").Value = "u"
$ws1.Range("C2").Value = 0.1
$ws1.Range("D2").Value = 1

$ws1.Range("B3").Value = "u"
$ws1.Range("C3").Value = 0.1
$ws1.Range("D3").Value = 1

$ws1.Range("B4").Value = "u"
$ws1.Range("C4").Value = 0.1
$ws1.Range("D4").Value = 1

# Rows 5, 9, 13 (alpha self-interaction terms): same "u" label, D becomes
# -0.004 instead of -0.04. These three pick up a distinguishable font
# variant in the source workbook, so nudge the font to force Excel to
# register a second font/style entry instead of reusing style 1.
$ws1.Range("B5").Value = "u"
$ws1.Range("B5").Font.Bold = $true
$ws1.Range("D5").Value = -0.004

$ws1.Range("B9").Value = "u"
$ws1.Range("B9").Font.Bold = $true
$ws1.Range("D9").Value = -0.004

$ws1.Range("B13").Value = "u"
$ws1.Range("B13").Font.Bold = $true
$ws1.Range("D13").Value = -0.004

# --- initial_conditions sheet ------------------------------------------

$ws2.Range("C2").Value = 20
$ws2.Range("D2").Value = 20

$ws2.Range("C3").Value = 20
$ws2.Range("D3").Value = 20

$ws2.Range("C4").Value = 20
$ws2.Range("D4").Value = 20

# --- view state: second sheet becomes the active / selected tab --------

[void]$ws1.Range("E4").Select()
$ws2.Activate()
[void]$ws2.Range("F4").Select()
